$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: "Det er instruktørenes opgave aktivere projekter, samt tilføje
# elever til det aktive projekter." paragraph -> insert "at " before
# "aktivere" and split "det aktive" into "de" + " aktive".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(18)
$xml1 = '<w:p ' + $wns + '><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Det er instruktørenes opgave </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">at </w:t></w:r>' + `
  '<w:r><w:t>aktive</w:t></w:r>' + `
  '<w:r><w:t>re projekter, samt</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> tilføje elever til de</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> aktive proj</w:t></w:r>' + `
  '<w:r><w:t>ekt</w:t></w:r>' + `
  '<w:r><w:t>er</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: "Vi besluttet på mødet ..." -> "Vi besluttede" split off into its
# own run, rest of the paragraph (incl. the lastRenderedPageBreak run and the
# trailing run) is left untouched.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(20)
$xml2 = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Textbody"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t>Vi besluttede</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> på mødet at det ville være sådan systemet skal fungere. DVS. alle de opgaver som på nuværende tidspunkt bliver brugt til at opfylde de målpinde hver elev, bliver ifølge af </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">dette system betragtet som projekter. De filer som som udgør opgaveformuleringen til målpinde opgaverne, vil i </w:t></w:r>' + `
  '<w:r><w:lastRenderedPageBreak/><w:t>dette system blive betragtet som kravspecifikation.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">  </w:t></w:r>' + `
  '</w:p>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: move the "Gennemgang af website" section (heading + 2 paragraphs)
# so it follows the "Database" section instead of trailing the document, add
# a new "Rettet af ..." line after the credits paragraph, and relocate the
# _GoBack bookmark to the new last paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(24)
$pEnd = $d.Paragraphs.Item($d.Paragraphs.Count)
$moveRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$xml3 = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Textbody"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t>Vi aftalte at der ville komme et dia</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">gram over databasen i dette referat,  det kommer </w:t></w:r>' + `
  '<w:r><w:t>som vedhæftet fil i mailen.</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading2"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t>Gennemgang af website</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Textbody"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t>Gennemgangen af websitet nåede vi desværre ikke under dette møde.</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Textbody"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
  '<w:r><w:t>Referatet er skrevet af Emil Schytte Bækgaard</w:t></w:r>' + `
  '<w:r><w:br/><w:t>Rettet af Tor Soya og Alexander Müllertz</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$moveRange.InsertXML($xml3)
